# Add a new "2021" column (O) to the "рус,англ" sheet, mirroring the
# existing 2020 column (N): copy its formatting into O4:O5, then set the
# new year label and data value, and finally move the selection as in
# the authored file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (style) from the 2020 column (N4:N5) into the new
# 2021 column (O4:O5) so the new cells keep the same look (borders,
# font, number format, etc.) as the rest of the year columns.
$ws.Range("N4:N5").Copy()
$ws.Range("O4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new year header and value.
$ws.Cells.Item(4, 15).Value = 2021
$ws.Cells.Item(5, 15).Value = 1.5020015556876996

# Match the author's final selection state.
$ws.Range("Q5").Select()
